# Insert two new rows at 117-118, shifting existing data down,
# and populate them with the new Frutilla price records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 rows before the current row 117 (old data shifts down to 119+)
$ws.Rows("117:118").Insert()

# New row 117: Especial, 135, 5500/6000/5815, 831
$ws.Range("A117").Value = 3
$ws.Range("B117").Value = "Femacal de La Calera"
$ws.Range("C117").Value = "Coquimbo"
$ws.Range("D117").Value = "2021-11-08"
$ws.Range("E117").Value = 5
$ws.Range("F117").Value = "Fruta"
$ws.Range("G117").Value = 100101
$ws.Range("H117").Value = "Berries"
$ws.Range("I117").Value = 100112025
$ws.Range("J117").Value = "Frutilla"
$ws.Range("K117").Value = "Sin especificar"
$ws.Range("L117").Value = "Especial"
$ws.Range("M117").Value = 135
$ws.Range("N117").Value = 5500
$ws.Range("O117").Value = 6000
$ws.Range("P117").Value = 5815
$ws.Range("Q117").Value = "$/bandeja 7 kilos"
$ws.Range("R117").Value = "Provincia de Melipilla"
$ws.Range("S117").Value = 831
$ws.Range("T117").Value = 7

# New row 118: Segunda, 87, 4000/4000/4000, 571
$ws.Range("A118").Value = 3
$ws.Range("B118").Value = "Femacal de La Calera"
$ws.Range("C118").Value = "Coquimbo"
$ws.Range("D118").Value = "2021-11-08"
$ws.Range("E118").Value = 5
$ws.Range("F118").Value = "Fruta"
$ws.Range("G118").Value = 100101
$ws.Range("H118").Value = "Berries"
$ws.Range("I118").Value = 100112025
$ws.Range("J118").Value = "Frutilla"
$ws.Range("K118").Value = "Sin especificar"
$ws.Range("L118").Value = "Segunda"
$ws.Range("M118").Value = 87
$ws.Range("N118").Value = 4000
$ws.Range("O118").Value = 4000
$ws.Range("P118").Value = 4000
$ws.Range("Q118").Value = "$/bandeja 7 kilos"
$ws.Range("R118").Value = "Provincia de Melipilla"
$ws.Range("S118").Value = 571
$ws.Range("T118").Value = 7

# Match the date formatting used by the rest of column D
$ws.Range("D117:D118").NumberFormat = "YYYY-MM-DD HH:MM:SS"
